# Generate Report for Handback
#
# The "342e139d-02d7-42cc-a8a1-f439910ea39d" handback row (row 7) on both the
# zh-cn and de-de status sheets now has a resolved "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime", plus an error detail
# explaining that the handed-back file version is stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c1d7a2257554c051caee2f089b5971c33b0f461/e2e/342e139d-02d7-42cc-a8a1-f439910ea39d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/633d98f91c408d2a16486fbf46a9e9f672508234/e2e/342e139d-02d7-42cc-a8a1-f439910ea39d.md."
$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/633d98f91c408d2a16486fbf46a9e9f672508234/e2e/342e139d-02d7-42cc-a8a1-f439910ea39d.md"
$latestMdDisplay = "342e139d-02d7-42cc-a8a1-f439910ea39d.md"

# ---------------------------------------------------------------------
# zh-cn sheet (sheet index 2), row 7
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item(2)

$wsZh.Range("I7").Value = $latestMdDisplay
$wsZh.Range("I7").Font.Underline = 2
$wsZh.Range("I7").Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestMdUrl, "", "", $latestMdDisplay) | Out-Null

$wsZh.Range("J7").Value = "342e139d-02d7-42cc-a8a1-f439910ea39d.cd3eaeccb604f139c16f2c86ce3149bd761c3225.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-24 20:57:30"
$wsZh.Range("P7").Value = $errorDetail

# ---------------------------------------------------------------------
# de-de sheet (sheet index 3), row 7
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item(3)

$wsDe.Range("I7").Value = $latestMdDisplay
$wsDe.Range("I7").Font.Underline = 2
$wsDe.Range("I7").Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestMdUrl, "", "", $latestMdDisplay) | Out-Null

$wsDe.Range("J7").Value = "342e139d-02d7-42cc-a8a1-f439910ea39d.cd3eaeccb604f139c16f2c86ce3149bd761c3225.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-24 20:57:38"
$wsDe.Range("P7").Value = $errorDetail
